$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row and data row to reflect the new "Metrics generation" layout
$ws.Range("B2").Value = "TestcasesbyModule"
$ws.Range("C1").Value = "Description"
$ws.Range("C2").Value = "Refering Modules Filename"

# Move the active selection to B4 (as recorded in the saved view state)
$ws.Range("B4").Select()
